$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "% of Q Drop's" column (column I) to the grade distribution table.
$ws.Range("I1").Value = "% of Q Drop's"

# Give the new data row the same "0.00%" value used by the other percentage
# columns on that row; copy it over so it reuses the existing shared string /
# style instead of being re-interpreted as a numeric percentage.
$ws.Range("E3").Copy()
$ws.Range("I3").PasteSpecial(-4163)
